# Applies the Efnb2-Epha3 sheet update (rows 2-13, full recompute incl. M2 target cluster)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Num($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value2 = [double]"$val"
}

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Efnb2"
$ws.Cells.Item(2, 3).Value2 = "Epha3"
$ws.Cells.Item(2, 4).Value2 = "FAPs"
Set-Num 2 5 "3"
Set-Num 2 6 "1"
Set-Num 2 7 "28.95628266666667"
Set-Num 2 8 "86.868848"
Set-Num 2 9 "0.5491054194301004"
Set-Num 2 10 "0.5491054194301005"
Set-Num 2 11 "3"
Set-Num 2 12 "1"
Set-Num 2 13 "46.25093466666667"
Set-Num 2 14 "138.752804"
Set-Num 2 15 "0.9569015955251317"
Set-Num 2 16 "0.9569015955251318"
Set-Num 2 17 "1339.255137805533"
Set-Num 2 18 "12053.29624024979"
Set-Num 2 19 "0.5254398519641597"
Set-Num 2 20 "0.5254398519641599"

# Row 3: ECs -> M2
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Efnb2"
$ws.Cells.Item(3, 3).Value2 = "Epha3"
$ws.Cells.Item(3, 4).Value2 = "M2"
Set-Num 3 5 "3"
Set-Num 3 6 "1"
Set-Num 3 7 "28.95628266666667"
Set-Num 3 8 "86.868848"
Set-Num 3 9 "0.5491054194301004"
Set-Num 3 10 "0.5491054194301005"
Set-Num 3 11 "2"
Set-Num 3 12 "0.6666666666666666"
Set-Num 3 13 "0.01632333333333333"
Set-Num 3 14 "0.04897"
Set-Num 3 15 "0.0003377190931065126"
Set-Num 3 16 "0.0003377190931065127"
Set-Num 3 17 "0.4726630540622222"
Set-Num 3 18 "4.25396748656"
Set-Num 3 19 "0.0001854433842698047"
Set-Num 3 20 "0.0001854433842698048"

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Efnb2"
$ws.Cells.Item(4, 3).Value2 = "Epha3"
$ws.Cells.Item(4, 4).Value2 = "sCs"
Set-Num 4 5 "3"
Set-Num 4 6 "1"
Set-Num 4 7 "28.95628266666667"
Set-Num 4 8 "86.868848"
Set-Num 4 9 "0.5491054194301004"
Set-Num 4 10 "0.5491054194301005"
Set-Num 4 11 "3"
Set-Num 4 12 "1"
Set-Num 4 13 "2.066797333333333"
Set-Num 4 14 "6.200391999999999"
Set-Num 4 15 "0.04276068538176181"
Set-Num 4 16 "0.04276068538176181"
Set-Num 4 17 "59.84676779871288"
Set-Num 4 18 "538.6209101884159"
Set-Num 4 19 "0.02348012408167088"
Set-Num 4 20 "0.02348012408167089"

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Efnb2"
$ws.Cells.Item(5, 3).Value2 = "Epha3"
$ws.Cells.Item(5, 4).Value2 = "FAPs"
Set-Num 5 5 "3"
Set-Num 5 6 "1"
Set-Num 5 7 "12.691493"
Set-Num 5 8 "38.074479"
Set-Num 5 9 "0.2406720388519202"
Set-Num 5 10 "0.2406720388519202"
Set-Num 5 11 "3"
Set-Num 5 12 "1"
Set-Num 5 13 "46.25093466666667"
Set-Num 5 14 "138.752804"
Set-Num 5 15 "0.9569015955251317"
Set-Num 5 16 "0.9569015955251318"
Set-Num 5 17 "586.9934135654573"
Set-Num 5 18 "5282.940722089115"
Set-Num 5 19 "0.2302994579756889"
Set-Num 5 20 "0.2302994579756889"

# Row 6: FAPs -> M2
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Efnb2"
$ws.Cells.Item(6, 3).Value2 = "Epha3"
$ws.Cells.Item(6, 4).Value2 = "M2"
Set-Num 6 5 "3"
Set-Num 6 6 "1"
Set-Num 6 7 "12.691493"
Set-Num 6 8 "38.074479"
Set-Num 6 9 "0.2406720388519202"
Set-Num 6 10 "0.2406720388519202"
Set-Num 6 11 "2"
Set-Num 6 12 "0.6666666666666666"
Set-Num 6 13 "0.01632333333333333"
Set-Num 6 14 "0.04897"
Set-Num 6 15 "0.0003377190931065126"
Set-Num 6 16 "0.0003377190931065127"
Set-Num 6 17 "0.2071674707366666"
Set-Num 6 18 "1.86450723663"
Set-Num 6 19 "8.127954269716586e-05"
Set-Num 6 20 "8.127954269716587e-05"

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Efnb2"
$ws.Cells.Item(7, 3).Value2 = "Epha3"
$ws.Cells.Item(7, 4).Value2 = "sCs"
Set-Num 7 5 "3"
Set-Num 7 6 "1"
Set-Num 7 7 "12.691493"
Set-Num 7 8 "38.074479"
Set-Num 7 9 "0.2406720388519202"
Set-Num 7 10 "0.2406720388519202"
Set-Num 7 11 "3"
Set-Num 7 12 "1"
Set-Num 7 13 "2.066797333333333"
Set-Num 7 14 "6.200391999999999"
Set-Num 7 15 "0.04276068538176181"
Set-Num 7 16 "0.04276068538176181"
Set-Num 7 17 "26.23074388841866"
Set-Num 7 18 "236.076694995768"
Set-Num 7 19 "0.01029130133353411"
Set-Num 7 20 "0.01029130133353411"

# Row 8: M2 -> FAPs
$ws.Cells.Item(8, 1).Value2 = "M2"
$ws.Cells.Item(8, 2).Value2 = "Efnb2"
$ws.Cells.Item(8, 3).Value2 = "Epha3"
$ws.Cells.Item(8, 4).Value2 = "FAPs"
Set-Num 8 5 "3"
Set-Num 8 6 "1"
Set-Num 8 7 "0.4888703333333334"
Set-Num 8 8 "1.466611"
Set-Num 8 9 "0.009270573592685367"
Set-Num 8 10 "0.009270573592685367"
Set-Num 8 11 "3"
Set-Num 8 12 "1"
Set-Num 8 13 "46.25093466666667"
Set-Num 8 14 "138.752804"
Set-Num 8 15 "0.9569015955251317"
Set-Num 8 16 "0.9569015955251318"
Set-Num 8 17 "22.61070984747155"
Set-Num 8 18 "203.496388627244"
Set-Num 8 19 "0.00887102666227378"
Set-Num 8 20 "0.00887102666227378"

# Row 9: M2 -> M2
$ws.Cells.Item(9, 1).Value2 = "M2"
$ws.Cells.Item(9, 2).Value2 = "Efnb2"
$ws.Cells.Item(9, 3).Value2 = "Epha3"
$ws.Cells.Item(9, 4).Value2 = "M2"
Set-Num 9 5 "3"
Set-Num 9 6 "1"
Set-Num 9 7 "0.4888703333333334"
Set-Num 9 8 "1.466611"
Set-Num 9 9 "0.009270573592685367"
Set-Num 9 10 "0.009270573592685367"
Set-Num 9 11 "2"
Set-Num 9 12 "0.6666666666666666"
Set-Num 9 13 "0.01632333333333333"
Set-Num 9 14 "0.04897"
Set-Num 9 15 "0.0003377190931065126"
Set-Num 9 16 "0.0003377190931065127"
Set-Num 9 17 "0.007979993407777778"
Set-Num 9 18 "0.07181994067"
Set-Num 9 19 "3.130849706298886e-06"
Set-Num 9 20 "3.130849706298887e-06"

# Row 10: M2 -> sCs
$ws.Cells.Item(10, 1).Value2 = "M2"
$ws.Cells.Item(10, 2).Value2 = "Efnb2"
$ws.Cells.Item(10, 3).Value2 = "Epha3"
$ws.Cells.Item(10, 4).Value2 = "sCs"
Set-Num 10 5 "3"
Set-Num 10 6 "1"
Set-Num 10 7 "0.4888703333333334"
Set-Num 10 8 "1.466611"
Set-Num 10 9 "0.009270573592685367"
Set-Num 10 10 "0.009270573592685367"
Set-Num 10 11 "3"
Set-Num 10 12 "1"
Set-Num 10 13 "2.066797333333333"
Set-Num 10 14 "6.200391999999999"
Set-Num 10 15 "0.04276068538176181"
Set-Num 10 16 "0.04276068538176181"
Set-Num 10 17 "1.010395901279111"
Set-Num 10 18 "9.093563111512"
Set-Num 10 19 "0.0003964160807052882"
Set-Num 10 20 "0.0003964160807052882"

# Row 11: sCs -> FAPs
$ws.Cells.Item(11, 1).Value2 = "sCs"
$ws.Cells.Item(11, 2).Value2 = "Efnb2"
$ws.Cells.Item(11, 3).Value2 = "Epha3"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
Set-Num 11 5 "3"
Set-Num 11 6 "1"
Set-Num 11 7 "10.59691233333333"
Set-Num 11 8 "31.790737"
Set-Num 11 9 "0.200951968125294"
Set-Num 11 10 "0.200951968125294"
Set-Num 11 11 "3"
Set-Num 11 12 "1"
Set-Num 11 13 "46.25093466666667"
Set-Num 11 14 "138.752804"
Set-Num 11 15 "0.9569015955251317"
Set-Num 11 16 "0.9569015955251318"
Set-Num 11 17 "490.1170999973942"
Set-Num 11 18 "4411.053899976548"
Set-Num 11 19 "0.1922912589230093"
Set-Num 11 20 "0.1922912589230093"

# Row 12: sCs -> M2
$ws.Cells.Item(12, 1).Value2 = "sCs"
$ws.Cells.Item(12, 2).Value2 = "Efnb2"
$ws.Cells.Item(12, 3).Value2 = "Epha3"
$ws.Cells.Item(12, 4).Value2 = "M2"
Set-Num 12 5 "3"
Set-Num 12 6 "1"
Set-Num 12 7 "10.59691233333333"
Set-Num 12 8 "31.790737"
Set-Num 12 9 "0.200951968125294"
Set-Num 12 10 "0.200951968125294"
Set-Num 12 11 "2"
Set-Num 12 12 "0.6666666666666666"
Set-Num 12 13 "0.01632333333333333"
Set-Num 12 14 "0.04897"
Set-Num 12 15 "0.0003377190931065126"
Set-Num 12 16 "0.0003377190931065127"
Set-Num 12 17 "0.1729769323211111"
Set-Num 12 18 "1.55679239089"
Set-Num 12 19 "6.786531643324314e-05"
Set-Num 12 20 "6.786531643324315e-05"

# Row 13: sCs -> sCs
$ws.Cells.Item(13, 1).Value2 = "sCs"
$ws.Cells.Item(13, 2).Value2 = "Efnb2"
$ws.Cells.Item(13, 3).Value2 = "Epha3"
$ws.Cells.Item(13, 4).Value2 = "sCs"
Set-Num 13 5 "3"
Set-Num 13 6 "1"
Set-Num 13 7 "10.59691233333333"
Set-Num 13 8 "31.790737"
Set-Num 13 9 "0.200951968125294"
Set-Num 13 10 "0.200951968125294"
Set-Num 13 11 "3"
Set-Num 13 12 "1"
Set-Num 13 13 "2.066797333333333"
Set-Num 13 14 "6.200391999999999"
Set-Num 13 15 "0.04276068538176181"
Set-Num 13 16 "0.04276068538176181"
Set-Num 13 17 "21.90167015210044"
Set-Num 13 18 "197.115031368904"
Set-Num 13 19 "0.008592843885851525"
Set-Num 13 20 "0.008592843885851525"

Write-Output "edit applied"